$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Range("D2") "68.267.16"
$ws.Range("E2").Value = "  +1.79%  "

# Row 3
Set-TextValue $ws.Range("D3") "3.898.56"
$ws.Range("E3").Value = "  +0.89%  "

# Row 4
$ws.Range("E4").Value = "  +0.03%  "

# Row 5
Set-TextValue $ws.Range("D5") "481.61"
$ws.Range("E5").Value = "  +1.54%  "

# Row 6
Set-TextValue $ws.Range("D6") "145.00"
$ws.Range("E6").Value = "  -0.15%  "

# Row 7
Set-TextValue $ws.Range("D7") "0.621"
$ws.Range("E7").Value = "  -1.66%  "

# Row 9
Set-TextValue $ws.Range("D9") "0.724"
$ws.Range("E9").Value = "  -2.98%  "

# Row 10
Set-TextValue $ws.Range("D10") "0.167"
$ws.Range("E10").Value = "  +7.80%  "

# Row 11
Set-TextValue $ws.Range("D11") "0.0000355"
$ws.Range("E11").Value = "  +14.01%  "

# Row 12
Set-TextValue $ws.Range("D12") "42.74"

# Row 13
Set-TextValue $ws.Range("D13") "10.62"
$ws.Range("E13").Value = "  +2.21%  "

# Row 14
Set-TextValue $ws.Range("D14") "4.518.95"
$ws.Range("E14").Value = "  +0.67%  "

# Row 15
Set-TextValue $ws.Range("D15") "14.66"
$ws.Range("E15").Value = "  -1.14%  "

# Row 16
Set-TextValue $ws.Range("D16") "3.938.68"
$ws.Range("E16").Value = "  +3.21%  "

# Row 17
$ws.Range("E17").Value = "  -0.40%  "

# Row 18
Set-TextValue $ws.Range("D18") "19.76"
$ws.Range("E18").Value = "  -1.82%  "

# Row 19
Set-TextValue $ws.Range("D19") "1.13"
$ws.Range("E19").Value = "  -3.06%  "

# Row 20
Set-TextValue $ws.Range("D20") "68.285.93"
$ws.Range("E20").Value = "  +1.48%  "

# Row 21
Set-TextValue $ws.Range("D21") "436.59"
$ws.Range("E21").Value = "  +0.99%  "

# Row 22
Set-TextValue $ws.Range("D22") "14.78"
$ws.Range("E22").Value = "  -0.87%  "

# Row 23
Set-TextValue $ws.Range("D23") "3.37"
$ws.Range("E23").Value = "  +0.89%  "

# Row 24
Set-TextValue $ws.Range("D24") "87.91"
$ws.Range("E24").Value = "  -0.81%  "

# Row 25
Set-TextValue $ws.Range("D25") "11.62"
$ws.Range("E25").Value = "  +16.72%  "

# Row 26
Set-TextValue $ws.Range("D26") "3.58"
$ws.Range("E26").Value = "  -0.89%  "

# Row 27
Set-TextValue $ws.Range("D27") "10.48"
$ws.Range("E27").Value = "  +4.62%  "

# Row 28
Set-TextValue $ws.Range("D28") "38.05"
$ws.Range("E28").Value = "  +0.02%  "

# Row 29
Set-TextValue $ws.Range("D29") "5.81"
$ws.Range("E29").Value = "  +4.47%  "

# Row 30
Set-TextValue $ws.Range("D30") "707.61"
$ws.Range("E30").Value = "  -3.35%  "

# Row 31
$ws.Range("B31").Value = "Hedera"
$ws.Range("C31").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue $ws.Range("D31") "0.131"
$ws.Range("E31").Value = "  -2.84%  "

# Row 32
$ws.Range("B32").Value = "Cosmos"
$ws.Range("C32").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue $ws.Range("D32") "13.37"
$ws.Range("E32").Value = "  -3.76%  "

# Row 33
Set-TextValue $ws.Range("D33") "2.87"
$ws.Range("E33").Value = "  +2.83%  "

# Row 34
Set-TextValue $ws.Range("D34") "0.0₃0921"
$ws.Range("E34").Value = "  +36.99%  "

# Row 35
Set-TextValue $ws.Range("D35") "41.80"
$ws.Range("E35").Value = "  -4.02%  "

# Row 36
Set-TextValue $ws.Range("D36") "59.54"
$ws.Range("E36").Value = "  +2.64%  "

# Row 37
$ws.Range("B37").Value = "NEARProtocol"
$ws.Range("C37").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue $ws.Range("D37") "5.69"
$ws.Range("E37").Value = "  +3.84%  "

# Row 38
$ws.Range("B38").Value = "Kaspa"
$ws.Range("C38").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue $ws.Range("D38") "0.150"
$ws.Range("E38").Value = "  -6.18%  "

# Row 39
Set-TextValue $ws.Range("D39") "0.999"
$ws.Range("E39").Value = "  -0.11%  "

# Row 40
Set-TextValue $ws.Range("D40") "0.0474"
$ws.Range("E40").Value = "  -2.10%  "

# Row 41
Set-TextValue $ws.Range("D41") "3.06"
$ws.Range("E41").Value = "  +10.50%  "

# Row 42
$ws.Range("B42").Value = "ThetaToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
Set-TextValue $ws.Range("D42") "3.02"
$ws.Range("E42").Value = "  +3.54%  "

# Row 43
$ws.Range("B43").Value = "Fetch.AI"
$ws.Range("C43").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextValue $ws.Range("D43") "2.73"
$ws.Range("E43").Value = "  +5.97%  "

# Row 44
Set-TextValue $ws.Range("D44") "0.341"
$ws.Range("E44").Value = "  -2.29%  "

# Row 45
Set-TextValue $ws.Range("D45") "0.142"
$ws.Range("E45").Value = "  -0.29%  "

# Row 46
Set-TextValue $ws.Range("D46") "0.999"
$ws.Range("E46").Value = "  -0.20%  "

# Row 47
Set-TextValue $ws.Range("D47") "3.42"
$ws.Range("E47").Value = "  -1.30%  "

# Row 48
Set-TextValue $ws.Range("D48") "2.14"
$ws.Range("E48").Value = "  -0.78%  "

# Row 49
Set-TextValue $ws.Range("D49") "146.21"
$ws.Range("E49").Value = "  +1.52%  "

# Row 50
Set-TextValue $ws.Range("D50") "3.14"
$ws.Range("E50").Value = "  -1.82%  "

# Row 51
Set-TextValue $ws.Range("D51") "2.85"
$ws.Range("E51").Value = "  -1.53%  "
